# === Edit script: apply diff changes to CasosColombia.xlsx ===
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix scattered cell corrections across existing rows ---

# Row 24: BP24 numeric 1 -> text "NaN"
$ws.Range("BP24").Value = "NaN"

# Row 34: BP34 text "NaN" -> numeric 1
$ws.Range("BP34").Value = 1

# Row 37: W37 text "NaN" -> numeric 4 ; DC37 numeric 1 -> text "NaN"
$ws.Range("W37").Value = 4
$ws.Range("DC37").Value = "NaN"

# Row 54: CG54 text "NaN" -> numeric 2
$ws.Range("CG54").Value = 2

# Row 73: X73 text "NaN" -> numeric 3
$ws.Range("X73").Value = 3

# Row 121: CG121 numeric 42 -> text "NaN"
$ws.Range("CG121").Value = "NaN"

# Row 128: W128 text "NaN" -> numeric 153
$ws.Range("W128").Value = 153

# Row 167: W167 text "NaN" -> numeric 655
$ws.Range("W167").Value = 655

# Rows 170-174: W column off-by-one increments
$ws.Range("W170").Value = 800
$ws.Range("W171").Value = 822
$ws.Range("W172").Value = 866
$ws.Range("W173").Value = 938
$ws.Range("W174").Value = 969

# --- Append new data row 185 ---

$ws.Range("A185").Value = 44079
$ws.Range("B185").Value = 658456
$ws.Range("C185").Value = 2712
$ws.Range("D185").Value = 87377
$ws.Range("E185").Value = 64965
$ws.Range("F185").Value = 223522
$ws.Range("G185").Value = 25831
$ws.Range("H185").Value = 4235
$ws.Range("I185").Value = 3330
$ws.Range("J185").Value = 6625
$ws.Range("K185").Value = 6141
$ws.Range("L185").Value = 12688
$ws.Range("M185").Value = 3782
$ws.Range("N185").Value = 20740
$ws.Range("O185").Value = 25381
$ws.Range("P185").Value = 5837
$ws.Range("Q185").Value = 5881
$ws.Range("R185").Value = 13038
$ws.Range("S185").Value = 9763
$ws.Range("T185").Value = 14998
$ws.Range("U185").Value = 12647
$ws.Range("V185").Value = 3152
$ws.Range("W185").Value = 1528
$ws.Range("X185").Value = 6902
$ws.Range("Y185").Value = 21224
$ws.Range("Z185").Value = 12327
$ws.Range("AA185").Value = 8166
$ws.Range("AB185").Value = 49559
$ws.Range("AC185").Value = 1230
$ws.Range("AD185").Value = 281
$ws.Range("AE185").Value = 415
$ws.Range("AF185").Value = 452
$ws.Range("AG185").Value = 165
$ws.Range("AH185").Value = 122
$ws.Range("AI185").Value = 339
$ws.Range("AJ185").Value = 1981
$ws.Range("AK185").Value = 3630
$ws.Range("AL185").Value = 36611
$ws.Range("AM185").Value = 7478
$ws.Range("AN185").Value = 2424
$ws.Range("AO185").Value = 38887
$ws.Range("AP185").Value = 983
$ws.Range("AQ185").Value = 20815
$ws.Range("AR185").Value = 1473
$ws.Range("AS185").Value = 8752
$ws.Range("AT185").Value = 1559
$ws.Range("AU185").Value = 1580
$ws.Range("AV185").Value = 5207
$ws.Range("AW185").Value = 1677
$ws.Range("AX185").Value = 950
$ws.Range("AY185").Value = 2480
$ws.Range("AZ185").Value = 2643
$ws.Range("BA185").Value = 50688
$ws.Range("BB185").Value = 12705
$ws.Range("BC185").Value = 3357
$ws.Range("BD185").Value = 8089
$ws.Range("BE185").Value = 4590
$ws.Range("BF185").Value = 280
$ws.Range("BG185").Value = 1413
$ws.Range("BH185").Value = 2624
$ws.Range("BI185").Value = 733
$ws.Range("BJ185").Value = 2053
$ws.Range("BK185").Value = 8743
$ws.Range("BL185").Value = 8728
$ws.Range("BM185").Value = 9012
$ws.Range("BN185").Value = 13936
$ws.Range("BO185").Value = 1891
$ws.Range("BP185").Value = 830
$ws.Range("BQ185").Value = 9175
$ws.Range("BR185").Value = 7818
$ws.Range("BS185").Value = 9255
$ws.Range("BT185").Value = 1751
$ws.Range("BU185").Value = 1647
$ws.Range("BV185").Value = 3622
$ws.Range("BW185").Value = 3762
$ws.Range("BX185").Value = 1094
$ws.Range("BY185").Value = 4913
$ws.Range("BZ185").Value = 2698
$ws.Range("CA185").Value = 1456
$ws.Range("CB185").Value = 769
$ws.Range("CC185").Value = 2350
$ws.Range("CD185").Value = 2028
$ws.Range("CE185").Value = 1466
$ws.Range("CF185").Value = 1071
$ws.Range("CG185").Value = 5373
$ws.Range("CH185").Value = 1604
$ws.Range("CI185").Value = 1214
$ws.Range("CJ185").Value = 1394
$ws.Range("CK185").Value = 1788
$ws.Range("CL185").Value = 1638
$ws.Range("CM185").Value = 1954
$ws.Range("CN185").Value = 1263
$ws.Range("CO185").Value = 1110
$ws.Range("CP185").Value = 1111
$ws.Range("CQ185").Value = 656
$ws.Range("CR185").Value = 3096
$ws.Range("CS185").Value = 1135
$ws.Range("CT185").Value = 823
$ws.Range("CU185").Value = 798
$ws.Range("CV185").Value = 1402
$ws.Range("CW185").Value = 1325
$ws.Range("CX185").Value = 668
$ws.Range("CY185").Value = 764
$ws.Range("CZ185").Value = 1013
$ws.Range("DA185").Value = 1269
$ws.Range("DB185").Value = 1101
$ws.Range("DC185").Value = 1231
$ws.Range("DD185").Value = 936
$ws.Range("DE185").Value = 318
$ws.Range("DF185").Value = 343
$ws.Range("DG185").Value = 718
$ws.Range("DH185").Value = 639
$ws.Range("DI185").Value = 419
$ws.Range("DJ185").Value = 534
$ws.Range("DK185").Value = 347
$ws.Range("DL185").Value = 620
$ws.Range("DM185").Value = 719
$ws.Range("DN185").Value = 517
$ws.Range("DO185").Value = 480
$ws.Range("DP185").Value = 372
$ws.Range("DQ185").Value = 516
$ws.Range("DR185").Value = 122359
$ws.Range("DS185").Value = 278669
$ws.Range("DT185").Value = 12033
$ws.Range("DU185").Value = 120434
$ws.Range("DV185").Value = 74480
$ws.Range("DW185").Value = 33920
$ws.Range("DX185").Value = 9942

# --- Update last-selected cell to match target view state ---
$ws.Range("DQ183").Select()

